# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g1.4")

# Updated values for existing rows (2010-2022), columns B (BR Deflacionado),
# C (NE Deflacionado) and D (SE Deflacionado) - automatic data refresh.
$data = @{
    2  = @(23257.86610669455, 11418.82435341108, 13560.43311809023)
    3  = @(25192.59314690959, 12245.94221027538, 14257.77027193148)
    4  = @(26979.84043437291, 13211.21379819802, 15697.30052998581)
    5  = @(27976.04665397691, 13739.79027459498, 16072.90969847004)
    6  = @(29912.95664259601, 14745.8907374447,  16785.95219769734)
    7  = @(31911.5976261855,  15975.15151436092, 17673.68733391674)
    8  = @(34224.49724617274, 17607.09315882314, 18604.68706838786)
    9  = @(35211.13679492509, 18274.86689908277, 19513.7812222432)
    10 = @(36646.18938775777, 19083.90609078256, 20594.91845447646)
    11 = @(37893.92756225695, 19561.56719637316, 20960.20827855906)
    12 = @(40040.12552822272, 20906.9003924803,  21326.36830874367)
    13 = @(44932.83670743866, 22962.50837336785, 23155.15872464016)
    14 = @(48783.70604975082, 24759.83037637985, 25299.05850862464)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Append new row 15 for year 2023
$ws.Cells.Item(15, 1).Value = 2023
$ws.Cells.Item(15, 2).Value = 51300.70579350938
$ws.Cells.Item(15, 3).Value = 26237.41536180414
$ws.Cells.Item(15, 4).Value = 26006.98661973922
